$d = $word.ActiveDocument
$r = $d.Content
$xml = $r.WordOpenXML

# shapetype anchorId
$old = "w14:anchorId=`"04EC68FC`" id=`"_x0000_t32`" coordsize=`"21600,21600`" o:spt=`"32`" o:oned=`"t`" path=`"m,l21600,21600e`" filled=`"f`""
$new = "w14:anchorId=`"60756FF4`" id=`"_x0000_t32`" coordsize=`"21600,21600`" o:spt=`"32`" o:oned=`"t`" path=`"m,l21600,21600e`" filled=`"f`""
if (-not $xml.Contains($old)) { throw "Pattern not found: shapetype anchorId" }
$xml = $xml.Replace($old, $new)

# abstractNum0 nsid/tmpl
$old = "<w:nsid w:val=`"0F9162FD`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"96E68B24`"/>"
$new = "<w:nsid w:val=`"04C47813`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"78BE8CB8`"/>"
if (-not $xml.Contains($old)) { throw "Pattern not found: abstractNum0 nsid/tmpl" }
$xml = $xml.Replace($old, $new)

# abstractNum1 full block
$old = "<w:abstractNum w:abstractNumId=`"1`" w15:restartNumberingAfterBreak=`"0`"><w:nsid w:val=`"22C049C0`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"0562D79A`"/><w:lvl w:ilvl=`"0`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"720`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"1`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"1440`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"2`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"2160`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"3`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"2880`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"4`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"3600`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"5`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"4320`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"6`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"5040`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"7`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"5760`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"8`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"6480`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl></w:abstractNum>"
$new = "<w:abstractNum w:abstractNumId=`"1`" w15:restartNumberingAfterBreak=`"0`"><w:nsid w:val=`"115F2363`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"76E22E68`"/><w:lvl w:ilvl=`"0`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"●`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"720`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"1`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"○`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"1440`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"2`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"■`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"2160`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"3`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"●`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"2880`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"4`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"○`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"3600`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"5`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"■`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"4320`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"6`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"●`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"5040`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"7`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"○`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"5760`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"8`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"■`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"6480`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl></w:abstractNum>"
if (-not $xml.Contains($old)) { throw "Pattern not found: abstractNum1 full block" }
$xml = $xml.Replace($old, $new)

# abstractNum2 nsid/tmpl
$old = "<w:nsid w:val=`"3C3D2341`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"50A08D24`"/>"
$new = "<w:nsid w:val=`"62D71AE8`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"F85C73CE`"/>"
if (-not $xml.Contains($old)) { throw "Pattern not found: abstractNum2 nsid/tmpl" }
$xml = $xml.Replace($old, $new)

# abstractNum3 nsid/tmpl
$old = "<w:nsid w:val=`"3F2E167A`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"F23EDB2E`"/>"
$new = "<w:nsid w:val=`"72EF0B11`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"66729E92`"/>"
if (-not $xml.Contains($old)) { throw "Pattern not found: abstractNum3 nsid/tmpl" }
$xml = $xml.Replace($old, $new)

# abstractNum4 full block
$old = "<w:abstractNum w:abstractNumId=`"4`" w15:restartNumberingAfterBreak=`"0`"><w:nsid w:val=`"6D743E48`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"42809B56`"/><w:lvl w:ilvl=`"0`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"●`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"720`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"1`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"○`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"1440`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"2`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"■`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"2160`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"3`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"●`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"2880`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"4`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"○`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"3600`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"5`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"■`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"4320`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"6`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"●`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"5040`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"7`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"○`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"5760`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"8`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"■`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"6480`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl></w:abstractNum>"
$new = "<w:abstractNum w:abstractNumId=`"4`" w15:restartNumberingAfterBreak=`"0`"><w:nsid w:val=`"7C1C1F28`"/><w:multiLevelType w:val=`"multilevel`"/><w:tmpl w:val=`"029C7494`"/><w:lvl w:ilvl=`"0`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"720`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"1`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"1440`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"2`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"2160`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"3`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"2880`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"4`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"3600`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"5`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"4320`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"6`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"5040`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"7`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"5760`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl><w:lvl w:ilvl=`"8`"><w:start w:val=`"1`"/><w:numFmt w:val=`"bullet`"/><w:lvlText w:val=`"-`"/><w:lvlJc w:val=`"left`"/><w:pPr><w:ind w:left=`"6480`" w:hanging=`"360`"/></w:pPr><w:rPr><w:u w:val=`"none`"/></w:rPr></w:lvl></w:abstractNum>"
if (-not $xml.Contains($old)) { throw "Pattern not found: abstractNum4 full block" }
$xml = $xml.Replace($old, $new)

# num->abstractNum mapping
$old = "<w:num w:numId=`"1`"><w:abstractNumId w:val=`"3`"/></w:num><w:num w:numId=`"2`"><w:abstractNumId w:val=`"0`"/></w:num><w:num w:numId=`"3`"><w:abstractNumId w:val=`"2`"/></w:num><w:num w:numId=`"4`"><w:abstractNumId w:val=`"1`"/></w:num><w:num w:numId=`"5`"><w:abstractNumId w:val=`"4`"/></w:num>"
$new = "<w:num w:numId=`"1`"><w:abstractNumId w:val=`"1`"/></w:num><w:num w:numId=`"2`"><w:abstractNumId w:val=`"2`"/></w:num><w:num w:numId=`"3`"><w:abstractNumId w:val=`"0`"/></w:num><w:num w:numId=`"4`"><w:abstractNumId w:val=`"4`"/></w:num><w:num w:numId=`"5`"><w:abstractNumId w:val=`"3`"/></w:num>"
if (-not $xml.Contains($old)) { throw "Pattern not found: num->abstractNum mapping" }
$xml = $xml.Replace($old, $new)

# BalloonText rsid
$old = "<w:unhideWhenUsed/><w:rsid w:val=`"008D3F57`"/><w:pPr><w:spacing w:line=`"240`" w:lineRule=`"auto`"/></w:pPr>"
$new = "<w:unhideWhenUsed/><w:rsid w:val=`"00B769FD`"/><w:pPr><w:spacing w:line=`"240`" w:lineRule=`"auto`"/></w:pPr>"
if (-not $xml.Contains($old)) { throw "Pattern not found: BalloonText rsid" }
$xml = $xml.Replace($old, $new)

# BalloonTextChar rsid
$old = "<w:semiHidden/><w:rsid w:val=`"008D3F57`"/><w:rPr><w:rFonts w:ascii=`"Segoe UI`" w:hAnsi=`"Segoe UI`" w:cs=`"Segoe UI`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>"
$new = "<w:semiHidden/><w:rsid w:val=`"00B769FD`"/><w:rPr><w:rFonts w:ascii=`"Segoe UI`" w:hAnsi=`"Segoe UI`" w:cs=`"Segoe UI`"/><w:sz w:val=`"18`"/><w:szCs w:val=`"18`"/></w:rPr>"
if (-not $xml.Contains($old)) { throw "Pattern not found: BalloonTextChar rsid" }
$xml = $xml.Replace($old, $new)

$r.WordOpenXML = $xml
Write-Output "Edit applied successfully"
